# Updated cryptos list on Tue Jun  4 13:36:06 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) values for each
# coin row to the latest scraped snapshot. Numeric-looking price strings are
# entered with a leading apostrophe so Excel stores them as text (matching
# the source data's text-formatted price column, e.g. "17.80" not "17.8").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "69.358.30"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").Value = "3.774.03"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'645.61"
$ws.Range("E5").Value = "  +1.64%  "

$ws.Range("D6").Value = "'166.02"
$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("D7").Value = "3.777.94"
$ws.Range("E7").Value = "  -1.45%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("E10").Value = "  -2.41%  "

$ws.Range("D11").Value = "'0.456"
$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").Value = "'6.89"
$ws.Range("E12").Value = "  +3.26%  "

$ws.Range("E13").Value = "  -5.02%  "

$ws.Range("D14").Value = "'34.85"
$ws.Range("E14").Value = "  -3.41%  "

$ws.Range("D15").Value = "4.414.35"

$ws.Range("D16").Value = "3.773.60"
$ws.Range("E16").Value = "  -2.40%  "

$ws.Range("D17").Value = "69.351.31"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "'17.80"
$ws.Range("E18").Value = "  -1.84%  "

$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("E20").Value = "  -2.04%  "

$ws.Range("D21").Value = "'468.19"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("E22").Value = "  -1.62%  "

$ws.Range("D23").Value = "'0.708"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("E24").Value = "  -5.80%  "

$ws.Range("D25").Value = "'81.88"
$ws.Range("E25").Value = "  -2.27%  "

$ws.Range("D26").Value = "'12.34"
$ws.Range("E26").Value = "  +2.61%  "

$ws.Range("D27").Value = "'10.39"
$ws.Range("E27").Value = "  +2.75%  "

$ws.Range("E28").Value = "  -3.74%  "

$ws.Range("D30").Value = "3.926.86"
$ws.Range("E30").Value = "  -1.43%  "

$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("E32").Value = "  +2.10%  "

$ws.Range("E33").Value = "  -2.23%  "

$ws.Range("D34").Value = "'28.66"
$ws.Range("E34").Value = "  -2.30%  "

$ws.Range("D35").Value = "'0.172"
$ws.Range("E35").Value = "  +14.61%  "

$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("D37").Value = "3.732.48"
$ws.Range("E37").Value = "  -1.20%  "

$ws.Range("D38").Value = "'8.85"
$ws.Range("E38").Value = "  -2.63%  "

$ws.Range("E39").Value = "  -2.42%  "

$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("D41").Value = "'3.26"
$ws.Range("E41").Value = "  -7.19%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("D43").Value = "'0.958"
$ws.Range("E43").Value = "  -2.54%  "

$ws.Range("D45").Value = "'45.26"

$ws.Range("E46").Value = "  +2.54%  "

$ws.Range("D47").Value = "'156.60"
$ws.Range("E47").Value = "  -0.57%  "

$ws.Range("D48").Value = "'47.47"
$ws.Range("E48").Value = "  +0.30%  "

$ws.Range("E49").Value = "  -1.77%  "

$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").Value = "'8.38"
$ws.Range("E51").Value = "  -1.09%  "

